$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = "Construa edifícios que forneçam funções especiais."
$ws.Range("C3").Value = "Casas, quartéis, instalações de pesquisa e muito mais podem ser construídas nesta área comum."
$ws.Range("C4").Value = "Instalações de madeira podem ser construídas aqui."
$ws.Range("C5").Value = "Instalações de mineração podem ser construídas aqui."
$ws.Range("C6").Value = "Instalações relacionadas à magia e refinarias de jade podem ser construídas nesta área carregada de energia arcana."
$ws.Columns.Item(3).AutoFit()
Write-Host ("ColumnWidth: " + $ws.Columns.Item(3).ColumnWidth)
